$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 13470231
$ws.Range("I64").Value = 4080105.5
$ws.Range("J64").Value = 26319876
$ws.Range("K64").Value = 4080105.5
$ws.Range("L64").Value = 26319876
$ws.Range("M64").Value = -4079857.5
$ws.Range("N64").Value = -26320372
$ws.Range("H67").Value = 13470231
$ws.Range("I67").Value = 4080105.5
$ws.Range("J67").Value = 26319876
$ws.Range("K67").Value = 4080105.5
$ws.Range("L67").Value = 26319876
$ws.Range("M67").Value = -4079247.5
$ws.Range("N67").Value = -26321592
$ws.Range("H70").Value = 5520.2
$ws.Range("J70").Value = 6325
$ws.Range("L70").Value = 18975
$ws.Range("N70").Value = -19515
$ws.Range("H73").Value = 5520.2
$ws.Range("J73").Value = 6325
$ws.Range("L73").Value = 18975
$ws.Range("N73").Value = -20847
$ws.Range("H94").Value = 12353527
$ws.Range("I94").Value = 15877393
$ws.Range("K94").Value = 15877393
$ws.Range("M94").Value = -15876942
$ws.Range("H96").Value = 3775.65
$ws.Range("I96").Value = 363.125
$ws.Range("J96").Value = 6050.6665
$ws.Range("K96").Value = 1089.375
$ws.Range("L96").Value = 18151.9995
$ws.Range("M96").Value = 283.625
$ws.Range("N96").Value = -20897.9995
$ws.Range("H132").Value = 9619
$ws.Range("I132").Value = 2592.7778
$ws.Range("K132").Value = 7778.3334
$ws.Range("M132").Value = -5248.3334

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6732
$ws.Range("I32").Value = 3919.6667
$ws.Range("J32").Value = 14402
$ws.Range("K32").Value = 3919.6667
$ws.Range("L32").Value = 14402
$ws.Range("M32").Value = -3632.6667
$ws.Range("N32").Value = -14976
$ws.Range("H74").Value = 7814207.5
$ws.Range("I74").Value = 12501192
$ws.Range("K74").Value = 12501192
$ws.Range("M74").Value = -12500318
$ws.Range("H77").Value = 7814207.5
$ws.Range("I77").Value = 12501192
$ws.Range("K77").Value = 62505960
$ws.Range("M77").Value = -62501592
$ws.Range("H88").Value = 675.375
$ws.Range("I88").Value = 667.1667
$ws.Range("J88").Value = 700
$ws.Range("K88").Value = 667.1667
$ws.Range("L88").Value = 700
$ws.Range("M88").Value = -261.1667
$ws.Range("N88").Value = -1512
$ws.Range("H91").Value = 675.375
$ws.Range("I91").Value = 667.1667
$ws.Range("J91").Value = 700
$ws.Range("K91").Value = 667.1667
$ws.Range("L91").Value = 700
$ws.Range("M91").Value = 736.8333
$ws.Range("N91").Value = -3508
$ws.Range("H97").Value = 615.3333
$ws.Range("I97").Value = 748.1818
$ws.Range("K97").Value = 748.1818
$ws.Range("M97").Value = -252.1818
$ws.Range("H102").Value = 335747.66
$ws.Range("I102").Value = 596665.4
$ws.Range("J102").Value = 2352.7778
$ws.Range("K102").Value = 596665.4
$ws.Range("L102").Value = 2352.7778
$ws.Range("M102").Value = -595043.4
$ws.Range("N102").Value = -5596.7778
$ws.Range("H132").Value = 23449.451
$ws.Range("I132").Value = 31275.422
$ws.Range("J132").Value = 11058.333
$ws.Range("K132").Value = 93826.266
$ws.Range("L132").Value = 33174.999
$ws.Range("M132").Value = -91296.266
$ws.Range("N132").Value = -38234.999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2693.8845
$ws.Range("I20").Value = 2140.3157
$ws.Range("K20").Value = 2140.3157
$ws.Range("M20").Value = -1893.3157
$ws.Range("H54").Value = 9995
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 9995
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 9995
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -10963
$ws.Range("H94").Value = 652998.3
$ws.Range("I94").Value = 721624.4399999999
$ws.Range("J94").Value = 1050
$ws.Range("K94").Value = 721624.4399999999
$ws.Range("L94").Value = 1050
$ws.Range("M94").Value = -721173.4399999999
$ws.Range("N94").Value = -1952
$ws.Range("H99").Value = 1042990.06
$ws.Range("I99").Value = 1603589.4
$ws.Range("J99").Value = 1877
$ws.Range("K99").Value = 1603589.4
$ws.Range("L99").Value = 1877
$ws.Range("M99").Value = -1602091.4
$ws.Range("N99").Value = -4873
$ws.Range("H107").Value = 1738.5
$ws.Range("I107").Value = 1752.7693
$ws.Range("J107").Value = 1701.4
$ws.Range("K107").Value = 1752.7693
$ws.Range("L107").Value = 1701.4
$ws.Range("M107").Value = 167.2307000000001
$ws.Range("N107").Value = -5541.4
$ws.Range("H134").Value = 5448
$ws.Range("I134").Value = 3172.25
$ws.Range("K134").Value = 9516.75
$ws.Range("M134").Value = -6981.75
$ws.Range("H140").Value = 30000
$ws.Range("J140").Value = 30000
$ws.Range("L140").Value = 30000
$ws.Range("N140").Value = -40360

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 248.6923
$ws.Range("I7").Value = 141.57143
$ws.Range("K7").Value = 141.57143
$ws.Range("M7").Value = -28.57142999999999
$ws.Range("H22").Value = 722.9167
$ws.Range("I22").Value = 262.5
$ws.Range("K22").Value = 262.5
$ws.Range("M22").Value = 87.5
$ws.Range("H60").Value = 233333
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 233333
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 233333
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -234355
$ws.Range("H132").Value = 16685734
$ws.Range("I132").Value = 20848798
$ws.Range("K132").Value = 62546394
$ws.Range("M132").Value = -62543864

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 285.8
$ws.Range("I86").Value = 316
$ws.Range("K86").Value = 948
$ws.Range("M86").Value = 238
$ws.Range("H89").Value = 285.8
$ws.Range("I89").Value = 316
$ws.Range("K89").Value = 2844
$ws.Range("M89").Value = 3084
$ws.Range("H131").Value = 9436582
$ws.Range("J131").Value = 6669407
$ws.Range("L131").Value = 20008221
$ws.Range("N131").Value = -20018301
$ws.Range("H141").Value = 13067.823
$ws.Range("I141").Value = 6716.8
$ws.Range("K141").Value = 20150.4
$ws.Range("M141").Value = -14970.4

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 362.14285
$ws.Range("J97").Value = 362.33334
$ws.Range("L97").Value = 362.33334
$ws.Range("N97").Value = -1354.33334
$ws.Range("H113").Value = 2237
$ws.Range("I113").Value = 2237
$ws.Range("K113").Value = 2237
$ws.Range("M113").Value = -67

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1016.1818
$ws.Range("I22").Value = 941
$ws.Range("K22").Value = 941
$ws.Range("M22").Value = -646
$ws.Range("H27").Value = 1016.1818
$ws.Range("I27").Value = 941
$ws.Range("K27").Value = 941
$ws.Range("M27").Value = -834
$ws.Range("H61").Value = 4299.6
$ws.Range("I61").Value = 4166.3335
$ws.Range("J61").Value = 4499.5
$ws.Range("K61").Value = 4166.3335
$ws.Range("L61").Value = 4499.5
$ws.Range("M61").Value = -3964.3335
$ws.Range("N61").Value = -4903.5
$ws.Range("H82").Value = 3907811
$ws.Range("I82").Value = 5210065
$ws.Range("J82").Value = 1049.5
$ws.Range("K82").Value = 5210065
$ws.Range("L82").Value = 1049.5
$ws.Range("M82").Value = -5209704
$ws.Range("N82").Value = -1771.5
$ws.Range("H85").Value = 3907811
$ws.Range("I85").Value = 5210065
$ws.Range("J85").Value = 1049.5
$ws.Range("K85").Value = 5210065
$ws.Range("L85").Value = 1049.5
$ws.Range("M85").Value = -5208817
$ws.Range("N85").Value = -3545.5
$ws.Range("H93").Value = 2180.182
$ws.Range("J93").Value = 4993.25
$ws.Range("L93").Value = 4993.25
$ws.Range("N93").Value = -7489.25
$ws.Range("H100").Value = 2324
$ws.Range("I100").Value = 1840.4445
$ws.Range("J100").Value = 4500
$ws.Range("K100").Value = 1840.4445
$ws.Range("L100").Value = 4500
$ws.Range("M100").Value = -1299.4445
$ws.Range("N100").Value = -5582
$ws.Range("H113").Value = 4299.6
$ws.Range("I113").Value = 4166.3335
$ws.Range("J113").Value = 4499.5
$ws.Range("K113").Value = 4166.3335
$ws.Range("L113").Value = 4499.5
$ws.Range("M113").Value = -1996.3335
$ws.Range("N113").Value = -8839.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H20").Value = 15000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 15000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 15000
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -15480
$ws.Range("H92").Value = 165161
$ws.Range("J92").Value = 165161
$ws.Range("L92").Value = 165161
$ws.Range("N92").Value = -170153
$ws.Range("H107").Value = 1987.68
$ws.Range("I107").Value = 2348.3684
$ws.Range("K107").Value = 7045.1052
$ws.Range("M107").Value = -5125.1052
$ws.Range("H126").Value = 1755.4
$ws.Range("I126").Value = 1194.5
$ws.Range("K126").Value = 3583.5
$ws.Range("M126").Value = -1113.5
